$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.164.34'
$ws.Range('E2').Value = '  +3.52%  '
$ws.Range('D3').Value = '1.576.73'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.92%  '
$ws.Range('D5').Value = '213.03'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  -0.85%  '
$ws.Range('D8').Value = '23.42'
$ws.Range('E8').Value = '  +6.36%  '
$ws.Range('D9').Value = '0.251'
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').Value = '0.0599'
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').Value = '0.0885'
$ws.Range('E11').Value = '  +2.00%  '
$ws.Range('D12').Value = '1.803.32'
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').Value = '1.571.70'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').Value = '28.116.14'
$ws.Range('E16').Value = '  +3.46%  '
$ws.Range('D17').Value = '63.63'
$ws.Range('E17').Value = '  +2.16%  '
$ws.Range('D18').Value = '230.06'
$ws.Range('E18').Value = '  +6.37%  '
$ws.Range('D19').Value = '0.0₃0706'
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('D20').Value = '7.46'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('E21').Value = '  -0.88%  '
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').Value = '9.35'
$ws.Range('E23').Value = '  +1.33%  '
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = '152.53'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').Value = '15.24'
$ws.Range('E26').Value = '  +0.99%  '
$ws.Range('D27').Value = '6.58'
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('D30').Value = '1.15'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').Value = '0.0474'
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('E32').Value = '  -0.79%  '
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('D34').Value = '1.417.24'
$ws.Range('E34').Value = '  -2.33%  '
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('D36').Value = '1.04'
$ws.Range('E36').Value = '  -5.56%  '
$ws.Range('E37').Value = '  -1.13%  '
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = '0.542'
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('D40').Value = '2.48'
$ws.Range('E40').Value = '  +5.69%  '
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('D43').Value = '5.64'
$ws.Range('E43').Value = '  -3.34%  '
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('E45').Value = '  +4.77%  '
$ws.Range('D46').Value = '63.76'
$ws.Range('E46').Value = '  -1.49%  '
$ws.Range('D47').Value = '1.715.50'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').Value = '87.13'
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('D49').Value = '0.0₆0107'
$ws.Range('E49').Value = '  +3.23%  '
$ws.Range('D50').Value = '0.0526'
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('D51').Value = '0.0946'
$ws.Range('E51').Value = '  -1.40%  '
